$wb = $excel.ActiveWorkbook

# Work on the "Bütçe" sheet
$ws = $wb.Worksheets.Item("Bütçe")

# Row 4 (Haftalık gider kalemleri) - fill in values
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = 6000
$ws.Range("D4").Value = 4000
$ws.Range("E4").Value = 4000
$ws.Range("F4").Value = 6000
$ws.Range("G4").Value = 4000
$ws.Range("H4").Value = 0

# Row 14 - fill in values
$ws.Range("B14").Value = 0
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("F14").Value = 0

# Restore selection on Roller sheet (kept as-is, but it is no longer the active tab)
$ws1 = $wb.Worksheets.Item("Roller")
$ws1.Activate()
$ws1.Range("D11").Select()

# Select a cell in Bütçe sheet, then make Bütçe the active sheet/tab (last = active)
$ws.Activate()
$ws.Range("E26").Select()
